$d = $word.ActiveDocument

# 1. Mark the "search.png" inline picture run as NoProof (adds <w:noProof/> to its rPr).
$searchShape = $d.InlineShapes.Item(2)
$searchShape.Range.NoProofing = 1

# 2. Remove the two empty paragraphs that follow the delinquency-buttons paragraph,
#    keeping the final empty paragraph before the section break.
$d.Paragraphs.Item(29).Range.Delete()
$d.Paragraphs.Item(29).Range.Delete()
